{"js": "// Recolor the shaded (\"striped\") table rows: cell shading that is\n// currently E5E5E5 (light gray) becomes B3B3B3 (darker gray).\nconst OLD_FILL = \"#E5E5E5\";\nconst NEW_FILL = \"#B3B3B3\";\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\n// Load every row for every table in the document.\nfor (const table of tables.items) {\n  table.rows.load(\"items\");\n}\nawait context.sync();\n\n// Load every cell for every row.\nconst rows = [];\nfor (const table of tables.items) {\n  for (const row of table.rows.items) {\n    row.cells.load(\"items\");\n    rows.push(row);\n  }\n}\nawait context.sync();\n\n// Load the current shading color of every cell.\nconst cells = [];\nfor (const row of rows) {\n  for (const cell of row.cells.items) {\n    cell.load(\"shadingColor\");\n    cells.push(cell);\n  }\n}\nawait context.sync();\n\n// Only cells that currently carry the old fill get updated to the new one.\nfor (const cell of cells) {\n  const current = cell.shadingColor;\n  if (current && current.toUpperCase() === OLD_FILL) {\n    cell.shadingColor = NEW_FILL;\n  }\n}\nawait context.sync();\n", "ps1": "# Recolor the shaded (\"striped\") table rows: cell shading that is\n# currently E5E5E5 (light gray) becomes B3B3B3 (darker gray).\n$oldColor = 0xE5E5E5\n$newColor = 0xB3B3B3\n\n$d = $word.ActiveDocument\n\nforeach ($t in $d.Tables) {\n    foreach ($row in $t.Rows) {\n        foreach ($cell in $row.Cells) {\n            if ($cell.Shading.BackgroundPatternColor -eq $oldColor) {\n                $cell.Shading.BackgroundPatternColor = $newColor\n            }\n        }\n    }\n}\n"}
